$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3").Value = 1.9
$ws.Range("R3").Value = 1.95
$ws.Range("G4").Value = 1.67
$ws.Range("H4").Value = 3.75
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.38
$ws.Range("K4").Value = 2.1
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("S4").Value = 2.15
$ws.Range("T4").Value = 1.67
$ws.Range("U4").Value = 3.15
$ws.Range("V4").Value = 1.35
$ws.Range("W4").Value = 4
$ws.Range("X4").Value = 1.22
$ws.Range("Y4").Value = 1.44
$ws.Range("Z4").Value = 2.63
$ws.Range("AA4").Value = 2.1
$ws.Range("AB4").Value = 1.67
$ws.Range("AC4").Value = 6
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 9
$ws.Range("AG4").Value = 15
$ws.Range("AH4").Value = 34
$ws.Range("AI4").Value = 8.5
$ws.Range("AJ4").Value = 7
$ws.Range("AK4").Value = 21
$ws.Range("AL4").Value = 67
$ws.Range("AM4").Value = 351
$ws.Range("AN4").Value = 11
$ws.Range("AO4").Value = 23
$ws.Range("AS4").Value = 51
$ws.Range("G5").Value = 4.3
$ws.Range("H5").Value = 3.05
$ws.Range("I5").Value = 1.88
$ws.Range("J5").Value = 4.5
$ws.Range("K5").Value = 2.07
$ws.Range("L5").Value = 2.42
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 2.8
$ws.Range("S5").Value = 1.98
$ws.Range("T5").Value = 1.65
$ws.Range("W5").Value = 3.2
$ws.Range("X5").Value = 1.26
$ws.Range("AA5").Value = 1.78
$ws.Range("AB5").Value = 1.83
$ws.Range("AC5").Value = 11.5
$ws.Range("AD5").Value = 25
$ws.Range("AF5").Value = 80
$ws.Range("AG5").Value = 45
$ws.Range("AH5").Value = 45
$ws.Range("AI5").Value = 8.5
$ws.Range("AJ5").Value = 6
$ws.Range("AK5").Value = 14
$ws.Range("AL5").Value = 70
$ws.Range("AM5").Value = 600
$ws.Range("AN5").Value = 6.4
$ws.Range("AO5").Value = 8.5
$ws.Range("AP5").Value = 8.25
$ws.Range("AR5").Value = 16
$ws.Range("AS5").Value = 28
